# daily auto push: 2026-02-25 10:06 UTC
# A new reading for 2026/02/25 (水, 時刻=16, ランキング=33) was appended to the
# log. In the sheet it slots in chronologically right after the existing
# 2026/02/25 rows (which end at row 855), pushing every following row down
# by one (old row 856 -> new row 857, ..., old row 897 -> new row 898).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 856:897 down to 857:898 to make room for the new entry.
$ws.Rows.Item(856).Insert()

# Force column A to text formatting before writing the date-like string so
# it is stored as literal text ("2026/02/25"), matching the rest of the
# column, instead of being auto-parsed into a date serial number. Resetting
# the style back to "Normal" afterwards drops the temporary text format so
# the cell ends up with the same (default) styling as its neighbours.
$ws.Range("A856").NumberFormat = "@"
$ws.Range("A856").Value = "2026/02/25"
$ws.Range("A856").Style = "Normal"

$ws.Range("B856").Value = "水"
$ws.Range("C856").Value = 16
$ws.Range("D856").Value = 33
